# Reporting Financial Calculations.xlsx - "updated financial calculations
# with latest work progress" commit.
#
# Real, content-level changes made by the author (everything else in the
# recorded diff - shared-string renumbering, cellXfs/mergeCell reordering,
# a shared-formula's cached "ref" span - is a mechanical side effect of
# Excel re-serialising the package and isn't an independent edit):
#
#   1) The three quarterly column headers on "Sheet1" were relabelled to
#      match the actual reporting periods (the financial reporting periods
#      run Dec-March / March-June / July-Sep / Sep-Dec, not the old
#      March-May / May-August / Aug-Nov labels).
#   2) Percent-complete progress on three milestones (rows 11, 13, 15) was
#      updated to reflect the latest work progress; every cash-flow total
#      that derives from those percentages recalculates accordingly.
#   3) The view was scrolled/zoomed to the area the author was last working
#      in (zoom 88%, selection on S19) instead of the old selection/scroll
#      position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Relabel the quarterly reporting-period headers (row 1) ---
$ws.Range("M1").Value = "March-June"
$ws.Range("P1").Value = "July-Sep"
$ws.Range("S1").Value = "Sep - Dec"

# --- Latest work-progress percentages feeding the cash-flow calculations ---
$ws.Range("L11").Value = 0.75
$ws.Range("L13").Value = 1
$ws.Range("L15").Value = 0.25

# --- Leave the view where the author was last working ---
$ws.Select()
$excel.ActiveWindow.Zoom = 88
$ws.Range("S19").Select()
